$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.NumberFormat = '@'
$c.Value = '26.638.88'
$c.Style = 'Normal'
$c = $ws.Range("E2")
$c.NumberFormat = '@'
$c.Value = '  -0.09%  '
$c.Style = 'Normal'
$c = $ws.Range("D3")
$c.NumberFormat = '@'
$c.Value = '1.591.17'
$c.Style = 'Normal'
$c = $ws.Range("E3")
$c.NumberFormat = '@'
$c.Value = '  -0.01%  '
$c.Style = 'Normal'
$c = $ws.Range("E4")
$c.NumberFormat = '@'
$c.Value = '  +0.26%  '
$c.Style = 'Normal'
$c = $ws.Range("D5")
$c.NumberFormat = '@'
$c.Value = '210.64'
$c.Style = 'Normal'
$c = $ws.Range("E5")
$c.NumberFormat = '@'
$c.Value = '  -0.24%  '
$c.Style = 'Normal'
$c = $ws.Range("D6")
$c.NumberFormat = '@'
$c.Value = '0.516'
$c.Style = 'Normal'
$c = $ws.Range("E6")
$c.NumberFormat = '@'
$c.Value = '  +1.26%  '
$c.Style = 'Normal'
$c = $ws.Range("E7")
$c.NumberFormat = '@'
$c.Value = '  +0.22%  '
$c.Style = 'Normal'
$c = $ws.Range("E8")
$c.NumberFormat = '@'
$c.Value = '  -0.38%  '
$c.Style = 'Normal'
$c = $ws.Range("D10")
$c.NumberFormat = '@'
$c.Value = '19.43'
$c.Style = 'Normal'
$c = $ws.Range("E10")
$c.NumberFormat = '@'
$c.Value = '  -1.22%  '
$c.Style = 'Normal'
$c = $ws.Range("E11")
$c.NumberFormat = '@'
$c.Value = '  +0.46%  '
$c.Style = 'Normal'
$c = $ws.Range("D12")
$c.NumberFormat = '@'
$c.Value = '1.815.12'
$c.Style = 'Normal'
$c = $ws.Range("E12")
$c.NumberFormat = '@'
$c.Value = '  +0.05%  '
$c.Style = 'Normal'
$c = $ws.Range("D13")
$c.NumberFormat = '@'
$c.Value = '1.585.09'
$c.Style = 'Normal'
$c = $ws.Range("E13")
$c.NumberFormat = '@'
$c.Value = '  -0.51%  '
$c.Style = 'Normal'
$c = $ws.Range("E14")
$c.NumberFormat = '@'
$c.Value = '  -0.70%  '
$c.Style = 'Normal'
$c = $ws.Range("D15")
$c.NumberFormat = '@'
$c.Value = '0.520'
$c.Style = 'Normal'
$c = $ws.Range("E15")
$c.NumberFormat = '@'
$c.Value = '  -1.44%  '
$c.Style = 'Normal'
$c = $ws.Range("D16")
$c.NumberFormat = '@'
$c.Value = '64.38'
$c.Style = 'Normal'
$c = $ws.Range("E16")
$c.NumberFormat = '@'
$c.Value = '  -0.55%  '
$c.Style = 'Normal'
$c = $ws.Range("D17")
$c.NumberFormat = '@'
$c.Value = '26.625.11'
$c.Style = 'Normal'
$c = $ws.Range("D18")
$c.NumberFormat = '@'
$c.Value = '0.0₃0727'
$c.Style = 'Normal'
$c = $ws.Range("E18")
$c.NumberFormat = '@'
$c.Value = '  -0.16%  '
$c.Style = 'Normal'
$c = $ws.Range("E19")
$c.NumberFormat = '@'
$c.Value = '  +0.24%  '
$c.Style = 'Normal'
$c = $ws.Range("D20")
$c.NumberFormat = '@'
$c.Value = '206.82'
$c.Style = 'Normal'
$c = $ws.Range("E20")
$c.NumberFormat = '@'
$c.Value = '  -0.58%  '
$c.Style = 'Normal'
$c = $ws.Range("D21")
$c.NumberFormat = '@'
$c.Value = '6.77'
$c.Style = 'Normal'
$c = $ws.Range("E21")
$c.NumberFormat = '@'
$c.Value = '  +0.45%  '
$c.Style = 'Normal'
$c = $ws.Range("E22")
$c.NumberFormat = '@'
$c.Value = '  -0.67%  '
$c.Style = 'Normal'
$c = $ws.Range("E23")
$c.NumberFormat = '@'
$c.Value = '  -2.96%  '
$c.Style = 'Normal'
$c = $ws.Range("D24")
$c.NumberFormat = '@'
$c.Value = '8.84'
$c.Style = 'Normal'
$c = $ws.Range("E24")
$c.NumberFormat = '@'
$c.Value = '  -0.62%  '
$c.Style = 'Normal'
$c = $ws.Range("D25")
$c.NumberFormat = '@'
$c.Value = '145.65'
$c.Style = 'Normal'
$c = $ws.Range("E25")
$c.NumberFormat = '@'
$c.Value = '  -1.04%  '
$c.Style = 'Normal'
$c = $ws.Range("E26")
$c.NumberFormat = '@'
$c.Value = '  +0.13%  '
$c.Style = 'Normal'
$c = $ws.Range("D27")
$c.NumberFormat = '@'
$c.Value = '7.18'
$c.Style = 'Normal'
$c = $ws.Range("E27")
$c.NumberFormat = '@'
$c.Value = '  -2.21%  '
$c.Style = 'Normal'
$c = $ws.Range("E28")
$c.NumberFormat = '@'
$c.Value = '  +0.08%  '
$c.Style = 'Normal'
$c = $ws.Range("D29")
$c.NumberFormat = '@'
$c.Value = '15.21'
$c.Style = 'Normal'
$c = $ws.Range("E29")
$c.NumberFormat = '@'
$c.Value = '  -0.33%  '
$c.Style = 'Normal'
$c = $ws.Range("D30")
$c.NumberFormat = '@'
$c.Value = '0.0503'
$c.Style = 'Normal'
$c = $ws.Range("E30")
$c.NumberFormat = '@'
$c.Value = '  -0.06%  '
$c.Style = 'Normal'
$c = $ws.Range("E31")
$c.NumberFormat = '@'
$c.Value = '  -0.10%  '
$c.Style = 'Normal'
$c = $ws.Range("D32")
$c.NumberFormat = '@'
$c.Value = '3.22'
$c.Style = 'Normal'
$c = $ws.Range("E32")
$c.NumberFormat = '@'
$c.Value = '  -1.04%  '
$c.Style = 'Normal'
$c = $ws.Range("D33")
$c.NumberFormat = '@'
$c.Value = '0.662'
$c.Style = 'Normal'
$c = $ws.Range("E33")
$c.NumberFormat = '@'
$c.Value = '  -0.37%  '
$c.Style = 'Normal'
$c = $ws.Range("E34")
$c.NumberFormat = '@'
$c.Value = '  +0.12%  '
$c.Style = 'Normal'
$c = $ws.Range("D35")
$c.NumberFormat = '@'
$c.Value = '1.278.76'
$c.Style = 'Normal'
$c = $ws.Range("E35")
$c.NumberFormat = '@'
$c.Value = '  -3.48%  '
$c.Style = 'Normal'
$c = $ws.Range("E36")
$c.NumberFormat = '@'
$c.Value = '  +2.08%  '
$c.Style = 'Normal'
$c = $ws.Range("E37")
$c.NumberFormat = '@'
$c.Value = '  -1.60%  '
$c.Style = 'Normal'
$c = $ws.Range("E38")
$c.NumberFormat = '@'
$c.Value = '  -0.36%  '
$c.Style = 'Normal'
$c = $ws.Range("E39")
$c.NumberFormat = '@'
$c.Value = '  +1.03%  '
$c.Style = 'Normal'
$c = $ws.Range("E40")
$c.NumberFormat = '@'
$c.Value = '  +0.23%  '
$c.Style = 'Normal'
$c = $ws.Range("D41")
$c.NumberFormat = '@'
$c.Value = '5.41'
$c.Style = 'Normal'
$c = $ws.Range("E41")
$c.NumberFormat = '@'
$c.Value = '  +0.46%  '
$c.Style = 'Normal'
$c = $ws.Range("E42")
$c.NumberFormat = '@'
$c.Value = '  +0.95%  '
$c.Style = 'Normal'
$c = $ws.Range("E43")
$c.NumberFormat = '@'
$c.Value = '  -0.07%  '
$c.Style = 'Normal'
$c = $ws.Range("D44")
$c.NumberFormat = '@'
$c.Value = '63.19'
$c.Style = 'Normal'
$c = $ws.Range("E44")
$c.NumberFormat = '@'
$c.Value = '  -0.66%  '
$c.Style = 'Normal'
$c = $ws.Range("B45")
$c.NumberFormat = '@'
$c.Value = 'WEMIXToken'
$c.Style = 'Normal'
$c = $ws.Range("C45")
$c.NumberFormat = '@'
$c.Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$c.Style = 'Normal'
$c = $ws.Range("D45")
$c.NumberFormat = '@'
$c.Value = '0.917'
$c.Style = 'Normal'
$c = $ws.Range("E45")
$c.NumberFormat = '@'
$c.Value = '  +9.62%  '
$c.Style = 'Normal'
$c = $ws.Range("B46")
$c.NumberFormat = '@'
$c.Value = 'RocketPoolETH'
$c.Style = 'Normal'
$c = $ws.Range("C46")
$c.NumberFormat = '@'
$c.Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$c.Style = 'Normal'
$c = $ws.Range("D46")
$c.NumberFormat = '@'
$c.Value = '1.727.44'
$c.Style = 'Normal'
$c = $ws.Range("E46")
$c.NumberFormat = '@'
$c.Value = '  +0.02%  '
$c.Style = 'Normal'
$c = $ws.Range("D47")
$c.NumberFormat = '@'
$c.Value = '89.75'
$c.Style = 'Normal'
$c = $ws.Range("E47")
$c.NumberFormat = '@'
$c.Value = '  -0.23%  '
$c.Style = 'Normal'
$c = $ws.Range("D48")
$c.NumberFormat = '@'
$c.Value = '1.60'
$c.Style = 'Normal'
$c = $ws.Range("E48")
$c.NumberFormat = '@'
$c.Value = '  -0.80%  '
$c.Style = 'Normal'
$c = $ws.Range("E49")
$c.NumberFormat = '@'
$c.Value = '  +3.08%  '
$c.Style = 'Normal'
$c = $ws.Range("D50")
$c.NumberFormat = '@'
$c.Value = '0.0506'
$c.Style = 'Normal'
$c = $ws.Range("E50")
$c.NumberFormat = '@'
$c.Value = '  -0.83%  '
$c.Style = 'Normal'
$c = $ws.Range("B51")
$c.NumberFormat = '@'
$c.Value = 'EnergySwap'
$c.Style = 'Normal'
$c = $ws.Range("C51")
$c.NumberFormat = '@'
$c.Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$c.Style = 'Normal'
$c = $ws.Range("D51")
$c.NumberFormat = '@'
$c.Value = '7.45'
$c.Style = 'Normal'
$c = $ws.Range("E51")
$c.NumberFormat = '@'
$c.Value = '  -0.90%  '
$c.Style = 'Normal'
